# Apply "v1475" performance-run update:
#  - Sheet "Sponza" (sheet2.xml): add new column T with v1475 results
#  - Sheet "ComplexMesh" (sheet3.xml): add new column S with v1475 results
#  - Update conditional formatting ranges to include the new columns
#  - Update selections / active sheet to match the authored view state

$wb = $excel.ActiveWorkbook

$wsPart = $wb.Worksheets.Item(1)   # PartOfSponza
$wsSponza = $wb.Worksheets.Item(2) # Sponza
$wsMesh = $wb.Worksheets.Item(3)   # ComplexMesh

# ---------------------------------------------------------------------------
# Sheet2 (Sponza): new column T ("v1475")
# ---------------------------------------------------------------------------

# Copy formatting from column S into column T first, so the new cells get the
# same styles (header/data/avg-var rows) as their neighbours.
$wsSponza.Range("S1:S16").Copy()
$wsSponza.Range("T1:T16").PasteSpecial(-4122)  # xlPasteFormats

$wsSponza.Range("T1").Value = "v1475"

$tData = @(4543, 4547, 4528, 4592, 4597, 4573, 4542, 4530, 4536, 4531)
for ($i = 0; $i -lt $tData.Length; $i++) {
    $wsSponza.Cells.Item(2 + $i, 20).Value = $tData[$i]
}

$wsSponza.Range("T12").Formula = "=AVERAGE(T2:T11)"
$wsSponza.Range("T13").Formula = "=_xlfn.VAR.S(T2:T11)"
$wsSponza.Range("T14").Formula = "=1-_xlfn.T.TEST(S2:S11,T2:T11,2,3)"
$wsSponza.Range("T15").Formula = "=S12/T12"
$wsSponza.Range("T16").Formula = "=B12/T12"

# Extend the conditional formatting that highlights the diff-accept rows.
$sponzaCF = $wsSponza.Range("B15:S16").FormatConditions
for ($i = 1; $i -le $sponzaCF.Count; $i++) {
    $sponzaCF.Item($i).ModifyAppliesToRange($wsSponza.Range("B15:T16"))
}

# ---------------------------------------------------------------------------
# Sheet3 (ComplexMesh): new column S ("v1475")
# ---------------------------------------------------------------------------

$wsMesh.Range("R1:R16").Copy()
$wsMesh.Range("S1:S16").PasteSpecial(-4122)  # xlPasteFormats

$wsMesh.Range("S1").Value = "v1475"

$sData = @(3738, 3728, 3689, 3687, 3704, 3722, 3718, 3702, 3726, 3711)
for ($i = 0; $i -lt $sData.Length; $i++) {
    $wsMesh.Cells.Item(2 + $i, 19).Value = $sData[$i]
}

$wsMesh.Range("S12").Formula = "=AVERAGE(S2:S11)"
$wsMesh.Range("S13").Formula = "=_xlfn.VAR.S(S2:S11)"
$wsMesh.Range("S14").Formula = "=1-_xlfn.T.TEST(R2:R11,S2:S11,2,3)"
$wsMesh.Range("S15").Formula = "=R12/S12"
$wsMesh.Range("S16").Formula = "=B12/S12"

$meshCF = $wsMesh.Range("B15:R16").FormatConditions
for ($i = 1; $i -le $meshCF.Count; $i++) {
    $meshCF.Item($i).ModifyAppliesToRange($wsMesh.Range("B15:S16"))
}

# ---------------------------------------------------------------------------
# View state: selections on each sheet, and ComplexMesh becomes the active tab
# ---------------------------------------------------------------------------

$wsPart.Range("D34").Select()
$wsSponza.Range("T14").Select()
$wsMesh.Range("S15").Select()
